# Auto-generated script to apply updated market/profit values
# per the commit 'chore: update Sheets via scheduled runner'
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1000
$ws.Range("J17").Value = 1000
$ws.Range("L17").Value = 3000
$ws.Range("N17").Value = -3336
$ws.Range("H18").Value = 503.64285
$ws.Range("I18").Value = 254.25
$ws.Range("K18").Value = 254.25
$ws.Range("M18").Value = 29.75
$ws.Range("H28").Value = 2400.9333
$ws.Range("I28").Value = 2359.818
$ws.Range("K28").Value = 2359.818
$ws.Range("M28").Value = -1874.818
$ws.Range("H41").Value = 1756.5
$ws.Range("I41").Value = 1500.125
$ws.Range("J41").Value = 2098.3333
$ws.Range("K41").Value = 1500.125
$ws.Range("L41").Value = 2098.3333
$ws.Range("M41").Value = -1060.125
$ws.Range("N41").Value = -2978.3333
$ws.Range("H64").Value = 3741.6667
$ws.Range("J64").Value = 3615
$ws.Range("L64").Value = 3615
$ws.Range("N64").Value = -4111
$ws.Range("H67").Value = 3741.6667
$ws.Range("J67").Value = 3615
$ws.Range("L67").Value = 3615
$ws.Range("N67").Value = -5331
$ws.Range("H137").Value = 2169.186
$ws.Range("I137").Value = 1730.375
$ws.Range("J137").Value = 2723.4736
$ws.Range("K137").Value = 5191.125
$ws.Range("L137").Value = 8170.4208
$ws.Range("M137").Value = -2641.125
$ws.Range("N137").Value = -13270.4208
$ws.Range("H138").Value = 2072.2737
$ws.Range("I138").Value = 1332.4286
$ws.Range("J138").Value = 2282.2297
$ws.Range("K138").Value = 3997.2858
$ws.Range("L138").Value = 6846.6891
$ws.Range("M138").Value = 1142.7142
$ws.Range("N138").Value = -17126.6891

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2969.2632
$ws.Range("I63").Value = 2431.6775
$ws.Range("J63").Value = 5350
$ws.Range("K63").Value = 2431.6775
$ws.Range("L63").Value = 5350
$ws.Range("M63").Value = -1745.6775
$ws.Range("N63").Value = -6722
$ws.Range("H66").Value = 2969.2632
$ws.Range("I66").Value = 2431.6775
$ws.Range("J66").Value = 5350
$ws.Range("K66").Value = 12158.3875
$ws.Range("L66").Value = 26750
$ws.Range("M66").Value = -8726.387499999999
$ws.Range("N66").Value = -33614
$ws.Range("H74").Value = 2700.7693
$ws.Range("I74").Value = 1787
$ws.Range("K74").Value = 1787
$ws.Range("M74").Value = -913
$ws.Range("H77").Value = 2700.7693
$ws.Range("I77").Value = 1787
$ws.Range("K77").Value = 8935
$ws.Range("M77").Value = -4567
$ws.Range("H132").Value = 2349.8845
$ws.Range("I132").Value = 1569.4736
$ws.Range("K132").Value = 4708.4208
$ws.Range("M132").Value = -2178.4208
$ws.Range("H133").Value = 37688.57
$ws.Range("J133").Value = 37688.57
$ws.Range("L133").Value = 37688.57
$ws.Range("N133").Value = -42748.57

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2842.5862
$ws.Range("J86").Value = 3567.6667
$ws.Range("L86").Value = 3567.6667
$ws.Range("N86").Value = -5813.6667
$ws.Range("H89").Value = 2842.5862
$ws.Range("J89").Value = 3567.6667
$ws.Range("L89").Value = 17838.3335
$ws.Range("N89").Value = -29070.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1887.3077
$ws.Range("I31").Value = 1847.4595
$ws.Range("J31").Value = 2624.5
$ws.Range("K31").Value = 1847.4595
$ws.Range("L31").Value = 2624.5
$ws.Range("M31").Value = -1552.4595
$ws.Range("N31").Value = -3214.5
$ws.Range("H34").Value = 1887.3077
$ws.Range("I34").Value = 1847.4595
$ws.Range("J34").Value = 2624.5
$ws.Range("K34").Value = 1847.4595
$ws.Range("L34").Value = 2624.5
$ws.Range("M34").Value = -1645.4595
$ws.Range("N34").Value = -3028.5
$ws.Range("H58").Value = 3802.4565
$ws.Range("I58").Value = 1041.25
$ws.Range("K58").Value = 1041.25
$ws.Range("M58").Value = -838.25
$ws.Range("H99").Value = 1826.6666
$ws.Range("I99").Value = 1712
$ws.Range("J99").Value = 2400
$ws.Range("K99").Value = 1712
$ws.Range("L99").Value = 2400
$ws.Range("M99").Value = -214
$ws.Range("N99").Value = -5396
$ws.Range("H122").Value = 1093
$ws.Range("I122").Value = 965.5
$ws.Range("K122").Value = 2896.5
$ws.Range("M122").Value = -446.5
$ws.Range("H126").Value = 1826.6666
$ws.Range("I126").Value = 1712
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 5136
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -2666
$ws.Range("N126").Value = -12140
$ws.Range("H132").Value = 2447.4443
$ws.Range("I132").Value = 1859.6364
$ws.Range("K132").Value = 5578.9092
$ws.Range("M132").Value = -3048.9092
$ws.Range("H136").Value = 3802.4565
$ws.Range("I136").Value = 1041.25
$ws.Range("K136").Value = 3123.75
$ws.Range("M136").Value = -573.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 3000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 9000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -10996
$ws.Range("H78").Value = 3000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 3000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 27000
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -36984
$ws.Range("H131").Value = 16154066
$ws.Range("J131").Value = 32181.604
$ws.Range("L131").Value = 96544.81200000001
$ws.Range("N131").Value = -106624.812

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 22504118
$ws.Range("I70").Value = 16670964
$ws.Range("K70").Value = 16670964
$ws.Range("M70").Value = -16670694
$ws.Range("H73").Value = 22504118
$ws.Range("I73").Value = 16670964
$ws.Range("K73").Value = 16670964
$ws.Range("M73").Value = -16670028
$ws.Range("H126").Value = 1869.35
$ws.Range("I126").Value = 1563.4286
$ws.Range("J126").Value = 2583.1667
$ws.Range("K126").Value = 4690.2858
$ws.Range("L126").Value = 7749.500100000001
$ws.Range("M126").Value = -2220.2858
$ws.Range("N126").Value = -12689.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1527.8422
$ws.Range("I136").Value = 1148.2667
$ws.Range("K136").Value = 3444.800099999999
$ws.Range("M136").Value = -894.8000999999995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 66608
$ws.Range("J20").Value = 66608
$ws.Range("L20").Value = 66608
$ws.Range("N20").Value = -67088
$ws.Range("H132").Value = 1471.7819
$ws.Range("I132").Value = 1145.8889
$ws.Range("J132").Value = 2938.3
$ws.Range("K132").Value = 3437.6667
$ws.Range("L132").Value = 8814.900000000001
$ws.Range("M132").Value = -907.6666999999998
$ws.Range("N132").Value = -13874.9
